# items.xlsx — add "Cód barras" column (product barcode) to the
# registration/bulk-upload template.
#
# feature | i-754 | Se agregó el input para ingresar el código de barras
# al formulario de registro de productos, tambien se agregó la columna
# para subir los productos de forma masiva desde el excel items.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the barcode column, right after "Fec. Vencimiento" (S1)
$ws.Range("T1").Value = "Cód barras"

# Sample barcode values for the two example rows
$ws.Range("T2").Value = 10000001
$ws.Range("T3").Value = 10000002

# An extra (empty) cell further down the sheet that carries a distinct,
# underlined-font style — mirrors the style added to the workbook.
$styledCell = $ws.Range("S8")
$styledCell.Font.Name = "Calibri"
$styledCell.Font.Size = 11
$styledCell.Font.Underline = $true
$styledCell.Font.Color = 0

# Move the selection/view to the newly active area of the sheet
[void]$ws.Range("S8").Select()
